$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 824 (existing rows 824.. shift down to 827..)
$ws.Range("A824:A826").EntireRow.Insert()

# Populate the 3 new rows with the new weekly price entries (constant columns
# A,B,C,E,F,G,H,I,J,K,Q,R,T copied from the surrounding rows in this sheet).
$newRows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44714, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Maduro",         200, 12000, 12000, 12000, "$/caja 20 kilos", "Ecuador", 600, 20),
    @(3, "Femacal de La Calera", "Coquimbo", 44714, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Pintón",         280, 13000, 13000, 13000, "$/caja 20 kilos", "Ecuador", 650, 20),
    @(3, "Femacal de La Calera", "Coquimbo", 44714, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Primera Pintón", 360, 14000, 15000, 14556, "$/caja 20 kilos", "Ecuador", 728, 20)
)

for ($i = 0; $i -lt 3; $i++) {
    $rowNum = 824 + $i
    $rowData = $newRows[$i]
    for ($j = 0; $j -lt 20; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

# The date column uses a date-time number format (style index 2 in this sheet);
# make sure the newly inserted date cells keep that format.
$ws.Range("D824:D826").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Inserted 3 rows and populated new data"
